# Weekly update: add a new "Fruta / hortaliza, semanal" report block for
# Femacal de La Calera - Kiwi. Insert 3 new rows (Especial / Primera /
# Segunda) at the top of the time series (row 1192) and push the rest of
# the table down, which also naturally carries the previously-last block
# (old rows 1240-1242) down to the new final rows 1243-1245.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 1192:1194, shifting existing data down.
$ws.Rows("1192:1194").Insert()

# Row 1192 - Especial
$ws.Range("A1192").Value = 3
$ws.Range("B1192").Value = "Femacal de La Calera"
$ws.Range("C1192").Value = "Coquimbo"
$ws.Range("D1192").Value = 45147
$ws.Range("E1192").Value = 5
$ws.Range("F1192").Value = "Fruta"
$ws.Range("G1192").Value = 100101
$ws.Range("H1192").Value = "Berries"
$ws.Range("I1192").Value = 100101007
$ws.Range("J1192").Value = "Kiwi"
$ws.Range("K1192").Value = "Hayward"
$ws.Range("L1192").Value = "Especial"
$ws.Range("M1192").Value = 60
$ws.Range("N1192").Value = 9000
$ws.Range("O1192").Value = 9000
$ws.Range("P1192").Value = 9000
$ws.Range("Q1192").Value = "$/bandeja 10 kilos"
$ws.Range("R1192").Value = "Provincia de Curicó"
$ws.Range("S1192").Value = 900
$ws.Range("T1192").Value = 10

# Row 1193 - Primera
$ws.Range("A1193").Value = 3
$ws.Range("B1193").Value = "Femacal de La Calera"
$ws.Range("C1193").Value = "Coquimbo"
$ws.Range("D1193").Value = 45147
$ws.Range("E1193").Value = 5
$ws.Range("F1193").Value = "Fruta"
$ws.Range("G1193").Value = 100101
$ws.Range("H1193").Value = "Berries"
$ws.Range("I1193").Value = 100101007
$ws.Range("J1193").Value = "Kiwi"
$ws.Range("K1193").Value = "Hayward"
$ws.Range("L1193").Value = "Primera"
$ws.Range("M1193").Value = 65
$ws.Range("N1193").Value = 8000
$ws.Range("O1193").Value = 8000
$ws.Range("P1193").Value = 8000
$ws.Range("Q1193").Value = "$/bandeja 10 kilos"
$ws.Range("R1193").Value = "Provincia de Curicó"
$ws.Range("S1193").Value = 800
$ws.Range("T1193").Value = 10

# Row 1194 - Segunda
$ws.Range("A1194").Value = 3
$ws.Range("B1194").Value = "Femacal de La Calera"
$ws.Range("C1194").Value = "Coquimbo"
$ws.Range("D1194").Value = 45147
$ws.Range("E1194").Value = 5
$ws.Range("F1194").Value = "Fruta"
$ws.Range("G1194").Value = 100101
$ws.Range("H1194").Value = "Berries"
$ws.Range("I1194").Value = 100101007
$ws.Range("J1194").Value = "Kiwi"
$ws.Range("K1194").Value = "Hayward"
$ws.Range("L1194").Value = "Segunda"
$ws.Range("M1194").Value = 58
$ws.Range("N1194").Value = 7000
$ws.Range("O1194").Value = 7000
$ws.Range("P1194").Value = 7000
$ws.Range("Q1194").Value = "$/bandeja 10 kilos"
$ws.Range("R1194").Value = "Provincia de Curicó"
$ws.Range("S1194").Value = 700
$ws.Range("T1194").Value = 10
